$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the forecast date header
$ws.Range("B1").Value = "24/03/2023"

# Update the hourly values
$ws.Range("B2").Value = 148
$ws.Range("B3").Value = 221
$ws.Range("B4").Value = 189
$ws.Range("B5").Value = 180
$ws.Range("B6").Value = 153
$ws.Range("B7").Value = 137
$ws.Range("B8").Value = 140
$ws.Range("B9").Value = 134
$ws.Range("B10").Value = 145
$ws.Range("B11").Value = 146
$ws.Range("B12").Value = 157
$ws.Range("B13").Value = 122
$ws.Range("B14").Value = 80
$ws.Range("B15").Value = 51
$ws.Range("B16").Value = 19
$ws.Range("B17").Value = 12
